$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (shifts old rows 13-23 down to 14-24),
# matching the new "Docentes responsaveis:" value row being split out
# from the label row above it.
$ws.Rows.Item(13).Insert()

# Fix up the cell values that need new content after the shift.
$ws.Range("B10").Value = 'Propiciar ao discente conhecimento básico dos elementos e das principais ferramentas do SIG, habilitando-os para a utilização das novas tecnologias de geoprocessamento com enfoque em estudos ambientais e ações de mitigação e correção. Desenvolver atividades práticas utilizando software de geoprocessamento (ArcGis, QGis, Idrisi, etc.). Orientar o desenvolvimento de projetos relacionados à aplicação de geotecnologias em estudos ambientais.'
$ws.Range("C10").Value = 'Propiciar ao discente conhecimento básico dos elementos e das principais ferramentas do SIG, habilitando-os para a utilização das novas tecnologias de geoprocessamento com enfoque em estudos ambientais e ações de mitigação e correção. Desenvolver atividades práticas utilizando software de geoprocessamento (ArcGis, QGis, Idrisi, etc.). Orientar o desenvolvimento de projetos relacionados à aplicação de geotecnologias em estudos ambientais.'
$ws.Range("B13").Value = '9146830 - Danúbia Caporusso Bargos'
$ws.Range("C13").Value = '9146830 - Danúbia Caporusso Bargos'
$ws.Range("B14").Value = 'Fundamentos e elementos essenciais de um SIG; Bases de Dados espaciais. Aquisição e gerenciamento de dados em SIG; Funções do SIG e análise de dados; Sensoriamento remoto.'
$ws.Range("C14").Value = 'Fundamentos e elementos essenciais de um SIG; Bases de Dados espaciais. Aquisição e gerenciamento de dados em SIG; Funções do SIG e análise de dados; Sensoriamento remoto.'
$ws.Range("B16").Value = 'Introdução ao SIG; Elementos essenciais de um SIG; Representações da paisagem em ambiente computacional: modelos vetorial e matricial; Coleta, importação e conversão de dados; Manutenção, criação e edição de bases de dados geográficos; Entrada e visualização de dados; Visualização, criação e edição e operações com shapefiles; Georreferenciamento; classificação de imagens de satélites; Elaboração e Layout de mapas temáticos.'
$ws.Range("C16").Value = 'Introdução ao SIG; Elementos essenciais de um SIG; Representações da paisagem em ambiente computacional: modelos vetorial e matricial; Coleta, importação e conversão de dados; Manutenção, criação e edição de bases de dados geográficos; Entrada e visualização de dados; Visualização, criação e edição e operações com shapefiles; Georreferenciamento; classificação de imagens de satélites; Elaboração e Layout de mapas temáticos.'
$ws.Range("B19").Value = 'Aulas teóricas e práticas, visitas técnicas e exercícios dirigidos. Avaliação baseada em provas, exercícios e trabalhos práticos e relatórios.'
$ws.Range("C19").Value = 'Aulas teóricas e práticas, visitas técnicas e exercícios dirigidos. Avaliação baseada em provas, exercícios e trabalhos práticos e relatórios.'
$ws.Range("B20").Value = 'Média ponderada das notas atribuídas às provas, exercícios e trabalhos práticos e relatórios.'
$ws.Range("C20").Value = 'Média ponderada das notas atribuídas às provas, exercícios e trabalhos práticos e relatórios.'
$ws.Range("B21").Value = 'Provas e/ou exercícios dirigidos.'
$ws.Range("C21").Value = 'Provas e/ou exercícios dirigidos.'
$ws.Range("B22").Value = 'BURROUGH, P. A. Principles of Geographical Information Systems - Spatial Information Systems and Geoestatistics, Oxford: Clarendon Press, 1998.BURROUGH, P. A.; MCDONNELL, R. A. Principles of Geographical Information Systems. Oxford University Press, 1998.CÂMARA, G. & MEDEIROS, J. S. GIS para Meio Ambiente. INPE. São José dos Campos, SP, 1998.CROSTA, A. P. Processamento Digital de Imagens de Sensoriamento Remoto. Campinas – SP, 1992.FLORENZANO, T. G. Imagens de Satélite para Estudos Ambientais. Oficina de textos. São Paulo, 2002.IBGE. Noções Básicas de Cartografia. Rio de Janeiro. Coleção Manuais Técnicos em Geociências, 1999.LONGLEY, P. A.; GOODCHILD, M. F.; MAGUIRE, D. J.; RHIND, D. W. Geographic Information Systems and Science. John Wiley & Sons, 2001.MIRANDA, J. I.; Fundamentos de Sistemas de Informações Geográficas. Brasília, Embrapa, 2005.MOREIRA, M. A. Fundamentos do Sensoriamento Remoto e Metodologias de Aplicação. São José dos Campos – SP – INPE, 2001.SILVA, A.B. Sistemas de Informações Geo-referenciadas. Editora da Unicamp. Campinas. 1999.SILVA, A. B; Sistemas de informações Geo-referenciadas: conceitos e fundamentos. Campinas: Editora da Unicamp, 2003.SILVA, J.X. Geoprocessamento para Análise Ambiental. Rio de Janeiro. 2001.Bibliografia complementar:CARVALHO, M. S.; PINA, M. F.; SANTOS, S. M.  Conceitos Básicos de Sistemas de Informação Geográfica e Cartografia Aplicados à Saúde. Rede Interagencial de Informações para a Saúde. Brasília. Ministério da Saúde, 2000.DENT, B. D.  Cartography Thematic Map Design. 5th Edition. WCB/McGraw-Hill, 1999.MATOS, J. Fundamentos da Informação Geográfica. Lisboa, Lidel, 2008.MORAES NOVO, E. M. L. Sensoriamento Remoto – Princípios e Aplicações. 2ªEdição. São Paulo, 1992.'
$ws.Range("C22").Value = 'BURROUGH, P. A. Principles of Geographical Information Systems - Spatial Information Systems and Geoestatistics, Oxford: Clarendon Press, 1998.BURROUGH, P. A.; MCDONNELL, R. A. Principles of Geographical Information Systems. Oxford University Press, 1998.CÂMARA, G. & MEDEIROS, J. S. GIS para Meio Ambiente. INPE. São José dos Campos, SP, 1998.CROSTA, A. P. Processamento Digital de Imagens de Sensoriamento Remoto. Campinas – SP, 1992.FLORENZANO, T. G. Imagens de Satélite para Estudos Ambientais. Oficina de textos. São Paulo, 2002.IBGE. Noções Básicas de Cartografia. Rio de Janeiro. Coleção Manuais Técnicos em Geociências, 1999.LONGLEY, P. A.; GOODCHILD, M. F.; MAGUIRE, D. J.; RHIND, D. W. Geographic Information Systems and Science. John Wiley & Sons, 2001.MIRANDA, J. I.; Fundamentos de Sistemas de Informações Geográficas. Brasília, Embrapa, 2005.MOREIRA, M. A. Fundamentos do Sensoriamento Remoto e Metodologias de Aplicação. São José dos Campos – SP – INPE, 2001.SILVA, A.B. Sistemas de Informações Geo-referenciadas. Editora da Unicamp. Campinas. 1999.SILVA, A. B; Sistemas de informações Geo-referenciadas: conceitos e fundamentos. Campinas: Editora da Unicamp, 2003.SILVA, J.X. Geoprocessamento para Análise Ambiental. Rio de Janeiro. 2001.Bibliografia complementar:CARVALHO, M. S.; PINA, M. F.; SANTOS, S. M.  Conceitos Básicos de Sistemas de Informação Geográfica e Cartografia Aplicados à Saúde. Rede Interagencial de Informações para a Saúde. Brasília. Ministério da Saúde, 2000.DENT, B. D.  Cartography Thematic Map Design. 5th Edition. WCB/McGraw-Hill, 1999.MATOS, J. Fundamentos da Informação Geográfica. Lisboa, Lidel, 2008.MORAES NOVO, E. M. L. Sensoriamento Remoto – Princípios e Aplicações. 2ªEdição. São Paulo, 1992.'

# Split the old merged column definition (A:B) into a standalone
# column A definition, matching the cleaned-up column layout.
$ws.Columns.Item(1).ColumnWidth = $ws.Columns.Item(1).ColumnWidth()
